$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks numeric need to be forced to Text so they stay as literal strings
# (matching the source data which stores everything as text), without leaving the cell locked
# to a non-default style afterward.
$textForceCells = @("D5", "D6", "D7", "D9", "D10", "D11", "D12", "D14", "D16", "D17", "D19", "D20", "D22", "D23", "D24", "D25", "D26", "D28", "D29", "D30", "D31", "D32", "D34", "D36", "D37", "D39", "D41", "D43", "D47", "D48", "D49", "D50")
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values
$ws.Range("D2").Value = '42.915.85'
$ws.Range("E2").Value = '  +0.53%  '
$ws.Range("D3").Value = '2.531.82'
$ws.Range("E3").Value = '  -0.52%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = '317.76'
$ws.Range("E5").Value = '  +1.52%  '
$ws.Range("D6").Value = '97.03'
$ws.Range("E6").Value = '  +1.66%  '
$ws.Range("D7").Value = '0.575'
$ws.Range("E7").Value = '  -0.70%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("D9").Value = '0.538'
$ws.Range("E9").Value = '  -0.18%  '
$ws.Range("D10").Value = '35.86'
$ws.Range("E10").Value = '  -1.89%  '
$ws.Range("D11").Value = '0.0818'
$ws.Range("E11").Value = '  +0.29%  '
$ws.Range("D12").Value = '7.53'
$ws.Range("E12").Value = '  -2.14%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '2.866.97'
$ws.Range("E13").Value = '  +13.07%  '
$ws.Range("B14").Value = 'TRON'
$ws.Range("C14").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D14").Value = '0.109'
$ws.Range("E14").Value = '  -4.43%  '
$ws.Range("B15").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C15").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D15").Value = '2.920.27'
$ws.Range("E15").Value = '  -0.50%  '
$ws.Range("D16").Value = '15.11'
$ws.Range("E16").Value = '  -4.06%  '
$ws.Range("D17").Value = '0.851'
$ws.Range("E17").Value = '  -2.26%  '
$ws.Range("D18").Value = '42.956.78'
$ws.Range("E18").Value = '  +0.57%  '
$ws.Range("D19").Value = '6.86'
$ws.Range("E19").Value = '  +2.95%  '
$ws.Range("D20").Value = '12.58'
$ws.Range("E20").Value = '  -5.07%  '
$ws.Range("D21").Value = '0.0₃0967'
$ws.Range("E21").Value = '  -0.44%  '
$ws.Range("D22").Value = '69.69'
$ws.Range("E22").Value = '  -2.01%  '
$ws.Range("D23").Value = '253.28'
$ws.Range("E23").Value = '  -1.61%  '
$ws.Range("D24").Value = '2.96'
$ws.Range("E24").Value = '  +0.15%  '
$ws.Range("D25").Value = '2.06'
$ws.Range("E25").Value = '  +1.25%  '
$ws.Range("D26").Value = '26.39'
$ws.Range("E26").Value = '  -4.55%  '
$ws.Range("E27").Value = '  +0.06%  '
$ws.Range("B28").Value = 'InjectiveProtocol'
$ws.Range("C28").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D28").Value = '41.54'
$ws.Range("E28").Value = '  +5.19%  '
$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").Value = '2.41'
$ws.Range("E29").Value = '  +1.72%  '
$ws.Range("D30").Value = '10.41'
$ws.Range("E30").Value = '  +3.07%  '
$ws.Range("D31").Value = '5.88'
$ws.Range("E31").Value = '  -1.88%  '
$ws.Range("D32").Value = '157.26'
$ws.Range("E32").Value = '  +0.79%  '
$ws.Range("E33").Value = '  -0.49%  '
$ws.Range("D34").Value = '19.34'
$ws.Range("E34").Value = '  -0.98%  '
$ws.Range("E35").Value = '  -1.47%  '
$ws.Range("D36").Value = '2.72'
$ws.Range("E36").Value = '  +3.55%  '
$ws.Range("D37").Value = '0.0796'
$ws.Range("E37").Value = '  +0.62%  '
$ws.Range("E38").Value = '  +1.42%  '
$ws.Range("D39").Value = '2.46'
$ws.Range("E39").Value = '  +8.10%  '
$ws.Range("E40").Value = '  -0.62%  '
$ws.Range("D41").Value = '21.79'
$ws.Range("E41").Value = '  -11.76%  '
$ws.Range("E42").Value = '  +0.50%  '
$ws.Range("D43").Value = '3.81'
$ws.Range("E44").Value = '  +0.23%  '
$ws.Range("E45").Value = '  -2.49%  '
$ws.Range("D46").Value = '1.996.33'
$ws.Range("E46").Value = '  -2.81%  '
$ws.Range("D47").Value = '9.15'
$ws.Range("E47").Value = '  +2.34%  '
$ws.Range("D48").Value = '84.66'
$ws.Range("E48").Value = '  -0.66%  '
$ws.Range("D49").Value = '105.77'
$ws.Range("E49").Value = '  +3.29%  '
$ws.Range("D50").Value = '75.24'
$ws.Range("E50").Value = '  +0.08%  '
$ws.Range("D51").Value = '2.774.79'
$ws.Range("E51").Value = '  -0.52%  '

# Restore default styling on the cells we temporarily reformatted as Text
foreach ($addr in $textForceCells) {
    $ws.Range($addr).Style = "Normal"
}
